# BSC-HGP - Assignment 01 - Evandro Gomez Quintino.docx
#
# The student number "2960774" run is removed from the cover-sheet table
# cell, leaving the paragraph empty. Word tracks the location of the most
# recent edit with the hidden "_GoBack" bookmark, so once that text is
# deleted, Word moves "_GoBack" from wherever it used to be (the very end
# of the document, after "I had to cancel my idea.") to the spot where the
# text was removed. Every other bookmark in the document (the TOC
# "_Toc..." markers) keeps its relative order but gets renumbered because
# the freshly (re)inserted "_GoBack" bookmark claims id 0.

$d = $word.ActiveDocument

# Locate the "2960774" run inside the "Student number:" table cell (find
# only - no replacement performed here).
$hit = $d.Content
$found = $hit.Find.Execute("2960774")

if ($found) {
    $editStart = $hit.Start
    $editEnd = $hit.End

    # Remove the stale "_GoBack" bookmark sitting at the end of the document
    # (left over from whatever was last edited before this change).
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }

    # Re-create "_GoBack" collapsed at the spot the text is about to be
    # removed from (this is what Word itself does on a real edit) while the
    # run still exists, then delete the student-number text itself.
    $goBackRange = $d.Range($editStart, $editStart)
    $d.Bookmarks.Add("_GoBack", $goBackRange)

    $numberRange = $d.Range($editStart, $editEnd)
    $numberRange.Text = ""
}
